# Adapt column header formatting to respective input file names (#7)
#
# - Header columns A1:J1 (suffix "_old")  -> suffix "_FV2210"
# - Header column  K1   ("diff")          -> unchanged
# - Header columns L1:U1 (suffix "_new")  -> suffix "_FV2304"
# - Freeze the header row (row 1)
# - Wrap the used range A1:U67 in an Excel Table ("Table1") with AutoFilter

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$baseNames = @(
    "Segmentname",
    "Segmentgruppe",
    "Segment",
    "Datenelement",
    "Segment ID",
    "Code",
    "Qualifier",
    "Beschreibung",
    "Bedingungsausdruck",
    "Bedingung"
)

# Columns 1-10 (A-J): "<name>_old" -> "<name>_FV2210"
for ($i = 0; $i -lt $baseNames.Count; $i++) {
    $col = $i + 1
    $ws.Cells.Item(1, $col).Value = "$($baseNames[$i])_FV2210"
}

# Column 11 (K): "diff" stays as-is.

# Columns 12-21 (L-U): "<name>_new" -> "<name>_FV2304"
for ($i = 0; $i -lt $baseNames.Count; $i++) {
    $col = 11 + $i + 1
    $ws.Cells.Item(1, $col).Value = "$($baseNames[$i])_FV2304"
}

# Freeze the top (header) row.
$ws.Activate()
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
$ws.Range("A1").Select()

# Turn the used range into a native Excel table with an AutoFilter.
$range = $ws.Range("A1:U67")
$table = $ws.ListObjects.Add(1, $range, $null, 1)
$table.Name = "Table1"
